$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5103.30962648207
$ws.Range("C2").Value = 4197.87911950057
$ws.Range("F2").Value = 0.767523084104255

$ws.Range("B3").Value = 1207.12767057941
$ws.Range("C3").Value = 2144.55814052456
$ws.Range("F3").Value = 12.7593343727145

$ws.Range("B4").Value = 1096.63641792581
$ws.Range("C4").Value = 2116.42533528749
$ws.Range("F4").Value = 11.2791754317365

$ws.Range("B5").Value = 5049.38273576318
$ws.Range("C5").Value = 4741.82030337591
$ws.Range("F5").Value = 136.145814733864

$ws.Range("B6").Value = 5046.48585762156
$ws.Range("C6").Value = 4954.35228016733
$ws.Range("F6").Value = 144.970987064407

$ws.Range("B7").Value = 5032.32142456955
$ws.Range("C7").Value = 5117.99497066575
$ws.Range("F7").Value = 160.127851712342

$ws.Range("B8").Value = 5230.43751372964
$ws.Range("C8").Value = 5421.57440301543
$ws.Range("F8").Value = 164.522157678575

$ws.Range("B9").Value = 5230.43751372964
$ws.Range("C9").Value = 4717.74732117159
$ws.Range("F9").Value = 135.196029268415

$ws.Range("B10").Value = 1198.04912219854
$ws.Range("C10").Value = 2426.29989811248
$ws.Range("F10").Value = 24.9444854130805

$ws.Range("B11").Value = 1070.29084929054
$ws.Range("C11").Value = 2540.72811260171
$ws.Range("F11").Value = 29.2607022212989

$ws.Range("B12").Value = 5457.44872542722
$ws.Range("C12").Value = 4959.11012323562
$ws.Range("F12").Value = 149.185755867017

$ws.Range("B13").Value = 5457.44872542722
$ws.Range("C13").Value = 5123.00640687772
$ws.Range("F13").Value = 156.014767685438

$ws.Range("B14").Value = 4260.05383201679
$ws.Range("C14").Value = 4483.87264664182
$ws.Range("F14").Value = 109.570333151043

$ws.Range("B15").Value = 4260.05383201679
$ws.Range("C15").Value = 4469.31398012108
$ws.Range("F15").Value = 108.963722046012

